$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old (stale / bad foreign-key) seed data ---
$ws.Range("A1:C4").ClearContents()

# --- Column C (asset "type" foreign key) - write grouped by category so that
#     shared strings land in the same order as the authoritative workbook ---
$ws.Range("C1").Value = "Solar Panel"
$ws.Range("C3").Value = "Solar Panel"
$ws.Range("C2").Value = "USmart Charger"
$ws.Range("C4").Value = "USmart Charger"

# --- Column B (asset name) ---
$ws.Range("B1").Value = "Backyard Solar Panel"
$ws.Range("B3").Value = "Roof Wind Turbine "
$ws.Range("B2").Value = "Blue Tesla Model X"
$ws.Range("B4").Value = "Red Tesla Model 3"

# --- Column H (preference / deadline) ---
$ws.Range("H1").Value = "No preferences"
$ws.Range("H3").Value = "No preferences"
$ws.Range("H2").Value = "Deadline: 07:00:00"
$ws.Range("H4").Value = "Deadline: 06:30:00"

# --- Column A (owner id) ---
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 2

# --- Column D ---
$ws.Range("D1").Value = 15
$ws.Range("D2").Value = 7.7
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 6.6

# --- Column E ---
$ws.Range("E1").Value = 55
$ws.Range("E2").Value = 17
$ws.Range("E3").Value = 9
$ws.Range("E4").Value = 23

# --- Column F ---
$ws.Range("F1").Value = 200
$ws.Range("F2").Value = 100
$ws.Range("F3").Value = 110
$ws.Range("F4").Value = 75

# --- Column G (boolean flag) ---
$ws.Range("G1").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("G4").Value = $false

# --- Columns I & J ("True" / "False" stored as literal TEXT, not booleans) ---
# Format the cells as Text first, enter the values as formulas that evaluate
# to the literal strings, then convert the formulas to plain values via
# copy / paste-special so they land in the sheet as shared-string text cells
# (matching the source workbook's CSV-import provenance) rather than as
# real boolean cells.
$ws.Range("I1:I4").NumberFormat = "@"
$ws.Range("J1:J4").NumberFormat = "@"
$ws.Range("I1:I4").Formula = "=T(""True"")"
$ws.Range("J1:J4").Formula = "=T(""False"")"
$ws.Range("I1:J4").Copy()
$ws.Range("I1:J4").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Column widths (best-fit sizing, matching the authoritative widths) ---
$ws.Columns.Item(1).ColumnWidth = 14.65
$ws.Columns.Item(2).ColumnWidth = 25
$ws.Columns.Item(8).ColumnWidth = 12.8

# --- Sheet view / selection ---
$null = $ws.Range("K4").Select()

# --- Page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1
